$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its original text formatting so that
# numeric-looking values (e.g. "254.06") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: 'Bitcoin'
$ws.Cells.Item(2, 4).Value = '98.533.40'
$ws.Cells.Item(2, 5).Value = '  -0.10%  '

# Row 3: 'Ethereum'
$ws.Cells.Item(3, 4).Value = '3.323.51'
$ws.Cells.Item(3, 5).Value = '  -0.30%  '

# Row 4: 'TetherUSD'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5: 'Solana'
$ws.Cells.Item(5, 4).Value = '254.06'
$ws.Cells.Item(5, 5).Value = '  -1.55%  '

# Row 6: 'BNB'
$ws.Cells.Item(6, 4).Value = '640.04'
$ws.Cells.Item(6, 5).Value = '  +0.68%  '

# Row 7: 'XRP'
$ws.Cells.Item(7, 4).Value = '1.53'
$ws.Cells.Item(7, 5).Value = '  +13.01%  '

# Row 8: 'Dogecoin'
$ws.Cells.Item(8, 4).Value = '0.427'
$ws.Cells.Item(8, 5).Value = '  +7.57%  '

# Row 9: 'Cardano'
$ws.Cells.Item(9, 4).Value = '1.07'
$ws.Cells.Item(9, 5).Value = '  +22.43%  '

# Row 10: 'USDC'
$ws.Cells.Item(10, 5).Value = '  +0.01%  '

# Row 11: 'LidoStakedEther'
$ws.Cells.Item(11, 4).Value = '3.316.12'
$ws.Cells.Item(11, 5).Value = '  -0.48%  '

# Row 12: 'Avalanche' -> 'TRON'
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).Value = '0.206'
$ws.Cells.Item(12, 5).Value = '  +2.93%  '

# Row 13: 'TRON' -> 'Avalanche'
$ws.Cells.Item(13, 2).Value = 'Avalanche'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(13, 4).Value = '42.98'
$ws.Cells.Item(13, 5).Value = '  +19.07%  '

# Row 14: 'WrappedBTC' -> 'ShibaInu'
$ws.Cells.Item(14, 2).Value = 'ShibaInu'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(14, 4).Value = '0.0000270'
$ws.Cells.Item(14, 5).Value = '  +8.14%  '

# Row 15: 'ShibaInu' -> 'WrappedBTC'
$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '98.202.15'
$ws.Cells.Item(15, 5).Value = '  -0.06%  '

# Row 16: 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(16, 4).Value = '3.947.77'
$ws.Cells.Item(16, 5).Value = '  +0.03%  '

# Row 17: 'Toncoin'
$ws.Cells.Item(17, 4).Value = '5.48'
$ws.Cells.Item(17, 5).Value = '  -1.18%  '

# Row 18: 'WrappedEther'
$ws.Cells.Item(18, 4).Value = '3.323.55'
$ws.Cells.Item(18, 5).Value = '  -0.57%  '

# Row 19: 'Polkadot'
$ws.Cells.Item(19, 4).Value = '6.94'
$ws.Cells.Item(19, 5).Value = '  +12.92%  '

# Row 20: 'Chainlink'
$ws.Cells.Item(20, 4).Value = '16.47'
$ws.Cells.Item(20, 5).Value = '  +9.69%  '

# Row 21: 'BitcoinCash'
$ws.Cells.Item(21, 4).Value = '537.41'
$ws.Cells.Item(21, 5).Value = '  +9.79%  '

# Row 22: 'SuiNetwork'
$ws.Cells.Item(22, 4).Value = '3.47'
$ws.Cells.Item(22, 5).Value = '  -2.47%  '

# Row 23: 'Uniswap'
$ws.Cells.Item(23, 4).Value = '10.10'
$ws.Cells.Item(23, 5).Value = '  +8.23%  '

# Row 24: 'Stellar'
$ws.Cells.Item(24, 4).Value = '0.421'
$ws.Cells.Item(24, 5).Value = '  +52.52%  '

# Row 25: 'PEPE'
$ws.Cells.Item(25, 4).Value = '0.0000201'
$ws.Cells.Item(25, 5).Value = '  -3.68%  '

# Row 26: 'NEARProtocol'
$ws.Cells.Item(26, 4).Value = '5.98'
$ws.Cells.Item(26, 5).Value = '  +4.03%  '

# Row 27: 'Litecoin'
$ws.Cells.Item(27, 4).Value = '95.72'
$ws.Cells.Item(27, 5).Value = '  +7.89%  '

# Row 28: 'Aptos'
$ws.Cells.Item(28, 4).Value = '12.42'
$ws.Cells.Item(28, 5).Value = '  +2.84%  '

# Row 29: 'Hedera' -> 'WrappedeETH'
$ws.Cells.Item(29, 2).Value = 'WrappedeETH'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(29, 4).Value = '3.507.47'
$ws.Cells.Item(29, 5).Value = '  -0.20%  '

# Row 30: 'WrappedeETH' -> 'Hedera'
$ws.Cells.Item(30, 2).Value = 'Hedera'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(30, 4).Value = '0.148'
$ws.Cells.Item(30, 5).Value = '  +19.10%  '

# Row 31: 'InternetComputer(DFINITY)' -> 'Dai'
$ws.Cells.Item(31, 2).Value = 'Dai'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(31, 4).Value = '0.998'
$ws.Cells.Item(31, 5).Value = '  -0.23%  '

# Row 32: 'Dai' -> 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '10.96'
$ws.Cells.Item(32, 5).Value = '  +15.87%  '

# Row 33: 'Cronos'
$ws.Cells.Item(33, 4).Value = '0.189'
$ws.Cells.Item(33, 5).Value = '  -3.89%  '

# Row 34: 'Binance-PegBSC-USD'
$ws.Cells.Item(34, 5).Value = '  +0.44%  '

# Row 35: 'EthereumClassic'
$ws.Cells.Item(35, 4).Value = '28.87'
$ws.Cells.Item(35, 5).Value = '  +4.58%  '

# Row 36: 'PolygonEcosystemToken'
$ws.Cells.Item(36, 4).Value = '0.513'
$ws.Cells.Item(36, 5).Value = '  +11.57%  '

# Row 37: 'PancakeSwap'
$ws.Cells.Item(37, 4).Value = '2.06'
$ws.Cells.Item(37, 5).Value = '  +3.11%  '

# Row 38: 'RenderToken'
$ws.Cells.Item(38, 4).Value = '7.52'
$ws.Cells.Item(38, 5).Value = '  +2.28%  '

# Row 39: 'Kaspa'
$ws.Cells.Item(39, 4).Value = '0.153'
$ws.Cells.Item(39, 5).Value = '  +3.03%  '

# Row 40: 'Bittensor'
$ws.Cells.Item(40, 4).Value = '513.80'
$ws.Cells.Item(40, 5).Value = '  +0.87%  '

# Row 41: 'WhiteBITCoin'
$ws.Cells.Item(41, 4).Value = '24.72'
$ws.Cells.Item(41, 5).Value = '  -0.05%  '

# Row 42: 'Fetch.AI' -> 'MantraDAO'
$ws.Cells.Item(42, 2).Value = 'MantraDAO'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(42, 4).Value = '3.90'
$ws.Cells.Item(42, 5).Value = '  +4.35%  '

# Row 43: 'MantraDAO' -> 'Fetch.AI'
$ws.Cells.Item(43, 2).Value = 'Fetch.AI'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(43, 4).Value = '1.28'
$ws.Cells.Item(43, 5).Value = '  +0.76%  '

# Row 44: 'ARBITRUM'
$ws.Cells.Item(44, 4).Value = '0.804'
$ws.Cells.Item(44, 5).Value = '  +5.31%  '

# Row 45: 'USDe'
$ws.Cells.Item(45, 5).Value = '  -0.01%  '

# Row 46: 'dogwifhat'
$ws.Cells.Item(46, 4).Value = '3.16'
$ws.Cells.Item(46, 5).Value = '  -1.85%  '

# Row 47: 'VeChain'
$ws.Cells.Item(47, 4).Value = '0.0384'
$ws.Cells.Item(47, 5).Value = '  +20.21%  '

# Row 48: 'Cosmos' -> 'Stacks'
$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).Value = '2.01'
$ws.Cells.Item(48, 5).Value = '  +3.24%  '

# Row 49: 'Stacks' -> 'Monero'
$ws.Cells.Item(49, 2).Value = 'Monero'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(49, 4).Value = '163.59'
$ws.Cells.Item(49, 5).Value = '  +1.75%  '

# Row 50: 'Filecoin' -> 'OKB'
$ws.Cells.Item(50, 2).Value = 'OKB'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(50, 4).Value = '49.80'
$ws.Cells.Item(50, 5).Value = '  +8.80%  '

# Row 51: 'Monero' -> 'Cosmos'
$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(51, 4).Value = '7.64'
$ws.Cells.Item(51, 5).Value = '  +15.97%  '
